$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Opleidingen in de dataset")
$ws.Name = "Educational Programmes"
$ws.Select()
$ws.Range("E18").Select()
